# Regenerate save_data to use K (strikeouts) instead of Strike# for column G.
# This updates the G2:G37 values on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 7
    3  = 10
    4  = 8
    5  = 8
    6  = 9
    7  = 9
    8  = 8
    9  = 5
    10 = 7
    11 = 8
    12 = 9
    13 = 4
    14 = 4
    15 = 7
    16 = 11
    17 = 7
    18 = 10
    19 = 11
    20 = 12
    21 = 12
    22 = 8
    23 = 7
    24 = 6
    25 = 5
    26 = 4
    27 = 6
    28 = 7
    29 = 5
    30 = 6
    31 = 5
    32 = 2
    33 = 10
    34 = 6
    35 = 5
    36 = 2
    37 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
